# updates slide title and readme contents
# Change-Id: I2054e601a0eabdda5af1ca3b2e19d38fff003c8b
#
# Slide 1 ("Title Slide"): the title placeholder ("Title 1") is widened /
# shifted left and its text is expanded from "Unit Testing in C#" to
# "An Introduction to Unit Testing in C#". PowerPoint's Normal Autofit then
# shrinks the (now longer) text to 90% to keep it inside the box.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)

# Reposition / resize the title placeholder.
# EMU -> points (1 pt = 12700 EMU):
#   off  x: 7999414 -> 7733792  (629.8751 -> 608.96 pt)
#   off  y: 1051551  (unchanged)
#   ext cx: 3565524 -> 3966464  (280.7499 -> 312.32 pt)
#   ext cy: 2384898  (unchanged)
$sh.Left  = 608.96
$sh.Width = 312.32

# Update the title text.
$sh.TextFrame.TextRange.Text = "An Introduction to Unit Testing in C#"

# The longer title no longer fits the box at full size, so Normal Autofit
# kicks in and shrinks the text to 90%.
$sh.TextFrame.AutoSize = 2
$sh.TextFrame.AutofitFontScale = 90000
